$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: append an extra blank paragraph right after $para, then delete the
# paragraph mark that separates them. Word merges the (still-empty) run that
# belonged to the throw-away paragraph onto the end of $para, which is the
# only reliable way to leave a genuinely empty trailing <w:r> (no <w:t> at
# all) behind a run of freshly inserted text.
# ---------------------------------------------------------------------------
function Add-TrailingEmptyRun($para) {
    $para.Range.InsertParagraphAfter() | Out-Null
    $endOfPara = $para.Range.End
    $mark = $d.Range($endOfPara - 1, $endOfPara)
    $mark.Delete() | Out-Null
}

# ===========================================================================
# 1) Existing trailing empty paragraph gets its first (empty) run filled in.
# ===========================================================================
$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertionPoint.InsertAfter("-Test Result: can’t seem to open it in Github... getting ahead of myself?")

# ===========================================================================
# 2) New paragraph: "-Inductive and Deductive Research: ..."
# ===========================================================================
$paraA = $d.Paragraphs.Last
$paraA.Range.InsertParagraphAfter() | Out-Null

$paraB = $d.Paragraphs.Last
$paraBStart = $paraB.Range.Start
$paraB.Range.InsertAfter("-Inductive and Deductive Research: Inductive is where you don’t start with a hypothesis but with observations, and then you make a hypothesis based on them. Deductive is the opposite- you make a hypothesis and then observe to prove (or disprove) a hypothesis.")

$uStart = $paraBStart + 1
$uEnd = $uStart + "Inductive and Deductive Research".Length
$d.Range($uStart, $uEnd).Font.Underline = 1

Add-TrailingEmptyRun $paraB

# ===========================================================================
# 3) New paragraph: "-Idiographic vs. Nomothetic: ..."
# ===========================================================================
$paraC = $d.Paragraphs.Last
$paraC.Range.InsertParagraphAfter() | Out-Null

$paraD = $d.Paragraphs.Last
$paraDStart = $paraD.Range.Start
$paraD.Range.InsertAfter("-Idiographic vs. Nomothetic: specific vs general (looking at a case study as opposed to testing a universal law)")

$uStart2 = $paraDStart + 1
$uEnd2 = $uStart2 + "Idiographic vs. Nomothetic".Length
$d.Range($uStart2, $uEnd2).Font.Underline = 1

Add-TrailingEmptyRun $paraD

# Paragraph-mark formatting and every run in this paragraph carry <w:u w:val="none"/>
# explicitly (as opposed to simply omitting the element), so set it back to
# "none" everywhere in the paragraph once the underlined run above exists.
$paraD.Range.Font.Underline = 0

# ===========================================================================
# 4) New paragraph (final paragraph of the document): "-Basic goal of the unit: ..."
# ===========================================================================
$paraE = $d.Paragraphs.Last
$paraE.Range.InsertParagraphAfter() | Out-Null

$paraF = $d.Paragraphs.Last
$paraFStart = $paraF.Range.Start
$part1 = "-Basic goal of the unit: come up with a problem that you will likely encounter in your future research and see what possible ways, using technology, you can overcome it (don’t "
$part2 = "actually "
$part3 = "have to do it- just propose it)"
$paraF.Range.InsertAfter($part1 + $part2 + $part3)

$italicStart = $paraFStart + $part1.Length
$italicEnd = $italicStart + $part2.Length
$d.Range($italicStart, $italicEnd).Font.Italic = 1

Add-TrailingEmptyRun $paraF
$paraF.Range.Font.Underline = 0
